$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.62068945992989
$ws.Range("C2").Value = 18.44807086640716
$ws.Range("D2").Value = 15.49919720907028
$ws.Range("E2").Value = 16.97211816502823
$ws.Range("G2").Value = 3.815075401916696
$ws.Range("I2").Value = 53.53622640527481
$ws.Range("J2").Value = 9.739729812649351
$ws.Range("K2").Value = 18.8075060075755
$ws.Range("B3").Value = 15.69083141230107
$ws.Range("C3").Value = 18.19091245392739
$ws.Range("D3").Value = 15.45990632205511
$ws.Range("E3").Value = 16.93466489763113
$ws.Range("G3").Value = 3.819641744356343
$ws.Range("I3").Value = 52.2121847377158
$ws.Range("J3").Value = 9.755754977490245
$ws.Range("K3").Value = 18.80905440865423
$ws.Range("B4").Value = 15.74179337013751
$ws.Range("C4").Value = 18.03497958851263
$ws.Range("D4").Value = 15.43947277559551
$ws.Range("E4").Value = 16.91572328343193
$ws.Range("G4").Value = 3.822582002164648
$ws.Range("I4").Value = 51.37963797142439
$ws.Range("J4").Value = 9.767163624891609
$ws.Range("K4").Value = 18.81658669937464
$ws.Range("B5").Value = 15.76452835703656
$ws.Range("C5").Value = 17.97199442045404
$ws.Range("D5").Value = 15.43207723706962
$ws.Range("E5").Value = 16.90902644637945
$ws.Range("G5").Value = 3.823814668618487
$ws.Range("I5").Value = 51.03575159112155
$ws.Range("J5").Value = 9.772206775247309
$ws.Range("K5").Value = 18.82130614211048
$ws.Range("B6").Value = 15.76842179295386
$ws.Range("C6").Value = 17.96157140660948
$ws.Range("D6").Value = 15.43090553441075
$ws.Range("E6").Value = 16.90797622098617
$ws.Range("G6").Value = 3.824021439652549
$ws.Range("I6").Value = 50.97837999158762
$ws.Range("J6").Value = 9.773067968735061
$ws.Range("K6").Value = 18.82218928461571
$ws.Range("B7").Value = 15.74209203917855
$ws.Range("C7").Value = 18.03412780106391
$ws.Range("D7").Value = 15.43936926253708
$ws.Range("E7").Value = 16.91562882691037
$ws.Range("G7").Value = 3.822598486481352
$ws.Range("I7").Value = 51.37501847433804
$ws.Range("J7").Value = 9.767230043877241
$ws.Range("K7").Value = 18.81664367427357
$ws.Range("B8").Value = 15.64322556626647
$ws.Range("C8").Value = 18.35903442450615
$ws.Range("D8").Value = 15.4848857354796
$ws.Range("E8").Value = 16.95836441313256
$ws.Range("G8").Value = 3.816621643177129
$ws.Range("I8").Value = 53.0839389111492
$ws.Range("J8").Value = 9.744929303187858
$ws.Range("K8").Value = 18.80667125764428
$ws.Range("B9").Value = 15.51273517221979
$ws.Range("C9").Value = 19.00886418321134
$ws.Range("D9").Value = 15.60325605605296
$ws.Range("E9").Value = 17.07419717179128
$ws.Range("G9").Value = 3.805976751331427
$ws.Range("I9").Value = 56.26726507193784
$ws.Range("J9").Value = 9.713671445488929
$ws.Range("K9").Value = 18.83950364545992
$ws.Range("B10").Value = 15.45646574007272
$ws.Range("C10").Value = 19.49015122788628
$ws.Range("D10").Value = 15.70768026741745
$ws.Range("E10").Value = 17.17857422509169
$ws.Range("G10").Value = 3.798801274727817
$ws.Range("I10").Value = 58.48837945178165
$ws.Range("J10").Value = 9.698342001938389
$ws.Range("K10").Value = 18.89570583297352
$ws.Range("B11").Value = 15.43965048201046
$ws.Range("C11").Value = 19.70912537892671
$ws.Range("D11").Value = 15.75889601959027
$ws.Range("E11").Value = 17.2301707506362
$ws.Range("G11").Value = 3.795674841579951
$ws.Range("I11").Value = 59.47042043236833
$ws.Range("J11").Value = 9.693032209360785
$ws.Range("K11").Value = 18.92824135224922
$ws.Range("B12").Value = 15.43455913759073
$ws.Range("C12").Value = 19.79198359556193
$ws.Range("D12").Value = 15.77881574409916
$ws.Range("E12").Value = 17.25029304949333
$ws.Range("G12").Value = 3.794510572225745
$ws.Range("I12").Value = 59.83800395372939
$ws.Range("J12").Value = 9.691261148476142
$ws.Range("K12").Value = 18.94156155901364
$ws.Range("B13").Value = 15.43559868946393
$ws.Range("C13").Value = 19.77414243295896
$ws.Range("D13").Value = 15.77450244944952
$ws.Range("E13").Value = 17.24593352211513
$ws.Range("G13").Value = 3.794760447415964
$ws.Range("I13").Value = 59.75903249189793
$ws.Range("J13").Value = 9.691631913927544
$ws.Range("K13").Value = 18.93864840506705
$ws.Range("B14").Value = 15.43920596438345
$ws.Range("C14").Value = 19.71594380784711
$ws.Range("D14").Value = 15.76052434815777
$ws.Range("E14").Value = 17.23181457209483
$ws.Range("G14").Value = 3.795578663693456
$ws.Range("I14").Value = 59.50074898250226
$ws.Range("J14").Value = 9.69288169814123
$ws.Range("K14").Value = 18.92931719679597
$ws.Range("B15").Value = 15.44158210266124
$ws.Range("C15").Value = 19.68028535237347
$ws.Range("D15").Value = 15.75203051956144
$ws.Range("E15").Value = 17.22324208810355
$ws.Range("G15").Value = 3.796082398075858
$ws.Range("I15").Value = 59.3419775565697
$ws.Range("J15").Value = 9.693678446955403
$ws.Range("K15").Value = 18.92373165863147
$ws.Range("B16").Value = 15.45774254625863
$ws.Range("C16").Value = 19.47583563515181
$ws.Range("D16").Value = 15.70440720257568
$ws.Range("E16").Value = 17.17528441735964
$ws.Range("G16").Value = 3.799008352072188
$ws.Range("I16").Value = 58.42361122646348
$ws.Range("J16").Value = 9.698722522526092
$ws.Range("K16").Value = 18.89371961883971
$ws.Range("B17").Value = 15.46991565824714
$ws.Range("C17").Value = 19.35037259837198
$ws.Range("D17").Value = 15.6761369943233
$ws.Range("E17").Value = 17.14691251753719
$ws.Range("G17").Value = 3.800838490176573
$ws.Range("I17").Value = 57.85280051917253
$ws.Range("J17").Value = 9.702243307970216
$ws.Range("K17").Value = 18.87709184167425
$ws.Range("B18").Value = 15.47774318957311
$ws.Range("C18").Value = 19.27821674395109
$ws.Range("D18").Value = 15.660226730031
$ws.Range("E18").Value = 17.13098159304655
$ws.Range("G18").Value = 3.801904110450415
$ws.Range("I18").Value = 57.52182809256835
$ws.Range("J18").Value = 9.704424939547925
$ws.Range("K18").Value = 18.86818416151935
$ws.Range("B19").Value = 15.48053487396531
$ws.Range("C19").Value = 19.25378932561111
$ws.Range("D19").Value = 15.65490014654833
$ws.Range("E19").Value = 17.12565449428308
$ws.Range("G19").Value = 3.802267144257793
$ws.Range("I19").Value = 57.40931705633777
$ws.Range("J19").Value = 9.705190480021479
$ws.Range("K19").Value = 18.86528092287064
$ws.Range("B20").Value = 15.46853423333499
$ws.Range("C20").Value = 19.36372810537557
$ws.Range("D20").Value = 15.67911023993459
$ws.Range("E20").Value = 17.14989266900694
$ws.Range("G20").Value = 3.800642327342259
$ws.Range("I20").Value = 57.91384093458307
$ws.Range("J20").Value = 9.701852306929945
$ws.Range("K20").Value = 18.87879399795017
$ws.Range("B21").Value = 15.43811168343516
$ws.Range("C21").Value = 19.73304038107006
$ws.Range("D21").Value = 15.7646158666362
$ws.Range("E21").Value = 17.23594587267248
$ws.Range("G21").Value = 3.795337801925462
$ws.Range("I21").Value = 59.57673128180143
$ws.Range("J21").Value = 9.692508099356974
$ws.Range("K21").Value = 18.93203089316607
$ws.Range("B22").Value = 15.42567291850795
$ws.Range("C22").Value = 19.97401278945157
$ws.Range("D22").Value = 15.82355686751295
$ws.Range("E22").Value = 17.29558526540625
$ws.Range("G22").Value = 3.791985407035636
$ws.Range("I22").Value = 60.63840995857055
$ws.Range("J22").Value = 9.687798095062915
$ws.Range("K22").Value = 18.97264925062665
$ws.Range("B23").Value = 15.43162649356003
$ws.Range("C23").Value = 19.84545919271981
$ws.Range("D23").Value = 15.79182212454638
$ws.Range("E23").Value = 17.2634464544252
$ws.Range("G23").Value = 3.7937642275455
$ws.Range("I23").Value = 60.07413685506896
$ws.Range("J23").Value = 9.69018396614282
$ws.Range("K23").Value = 18.95043864439678
$ws.Range("B24").Value = 15.46915619427085
$ws.Range("C24").Value = 19.35769015489709
$ws.Range("D24").Value = 15.67776496766912
$ws.Range("E24").Value = 17.14854415689469
$ws.Range("G24").Value = 3.800730970655362
$ws.Range("I24").Value = 57.88625328540638
$ws.Range("J24").Value = 9.702028588000845
$ws.Range("K24").Value = 18.87802242256383
$ws.Range("B25").Value = 15.54114737017377
$ws.Range("C25").Value = 18.83211111636197
$ws.Range("D25").Value = 15.56814340613381
$ws.Range("E25").Value = 17.03945459804187
$ws.Range("G25").Value = 3.808742395717915
$ws.Range("I25").Value = 55.42572228450435
$ws.Range("J25").Value = 9.720788706449493
$ws.Range("K25").Value = 18.82499285460029
